$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.803092333333333
$ws.Range("H2").Value = 5.409276999999999
$ws.Range("I2").Value = 0.1744886524959502
$ws.Range("J2").Value = 0.1744886524959502
$ws.Range("M2").Value = 16.27546433333333
$ws.Range("N2").Value = 48.826393
$ws.Range("O2").Value = 0.06628560529319844
$ws.Range("P2").Value = 0.06628560529319844
$ws.Range("Q2").Value = 29.34616496087344
$ws.Range("R2").Value = 264.115484647861
$ws.Range("S2").Value = 0.01156608594748862
$ws.Range("T2").Value = 0.01156608594748862
# Row 3
$ws.Range("G3").Value = 1.803092333333333
$ws.Range("H3").Value = 5.409276999999999
$ws.Range("I3").Value = 0.1744886524959502
$ws.Range("J3").Value = 0.1744886524959502
$ws.Range("O3").Value = 0.3480686258826592
$ws.Range("P3").Value = 0.3480686258826592
$ws.Range("Q3").Value = 154.0980016351328
$ws.Range("R3").Value = 1386.882014716195
$ws.Range("S3").Value = 0.06073402550638223
$ws.Range("T3").Value = 0.06073402550638222
# Row 4
$ws.Range("G4").Value = 1.803092333333333
$ws.Range("H4").Value = 5.409276999999999
$ws.Range("I4").Value = 0.1744886524959502
$ws.Range("J4").Value = 0.1744886524959502
$ws.Range("M4").Value = 42.61351133333333
$ws.Range("N4").Value = 127.840534
$ws.Range("O4").Value = 0.17355341356458
$ws.Range("P4").Value = 0.17355341356458
$ws.Range("Q4").Value = 76.83609558154643
$ws.Range("R4").Value = 691.5248602339179
$ws.Range("S4").Value = 0.03028310126895593
$ws.Range("T4").Value = 0.03028310126895593
# Row 5
$ws.Range("G5").Value = 1.803092333333333
$ws.Range("H5").Value = 5.409276999999999
$ws.Range("I5").Value = 0.1744886524959502
$ws.Range("J5").Value = 0.1744886524959502
$ws.Range("M5").Value = 101.183272
$ws.Range("N5").Value = 303.549816
$ws.Range("O5").Value = 0.4120923552595624
$ws.Range("P5").Value = 0.4120923552595624
$ws.Range("Q5").Value = 182.4427820047813
$ws.Range("R5").Value = 1641.985038043032
$ws.Range("S5").Value = 0.07190543977312346
$ws.Range("T5").Value = 0.07190543977312344
# Row 6
$ws.Range("I6").Value = 0.4384883998568034
$ws.Range("J6").Value = 0.4384883998568034
$ws.Range("M6").Value = 16.27546433333333
$ws.Range("N6").Value = 48.826393
$ws.Range("O6").Value = 0.06628560529319844
$ws.Range("P6").Value = 0.06628560529319844
$ws.Range("Q6").Value = 73.74664616615021
$ws.Range("R6").Value = 663.7198154953519
$ws.Range("S6").Value = 0.02906546899855424
$ws.Range("T6").Value = 0.02906546899855424
# Row 7
$ws.Range("I7").Value = 0.4384883998568034
$ws.Range("J7").Value = 0.4384883998568034
$ws.Range("O7").Value = 0.3480686258826592
$ws.Range("P7").Value = 0.3480686258826592
$ws.Range("R7").Value = 3485.22191399924
$ws.Range("S7").Value = 0.1526240548036436
$ws.Range("T7").Value = 0.1526240548036436
# Row 8
$ws.Range("I8").Value = 0.4384883998568034
$ws.Range("J8").Value = 0.4384883998568034
$ws.Range("M8").Value = 42.61351133333333
$ws.Range("N8").Value = 127.840534
$ws.Range("O8").Value = 0.17355341356458
$ws.Range("P8").Value = 0.17355341356458
$ws.Range("Q8").Value = 193.0884107410862
$ws.Range("R8").Value = 1737.795696669776
$ws.Range("S8").Value = 0.07610115860361873
$ws.Range("T8").Value = 0.07610115860361873
# Row 9
$ws.Range("I9").Value = 0.4384883998568034
$ws.Range("J9").Value = 0.4384883998568034
$ws.Range("M9").Value = 101.183272
$ws.Range("N9").Value = 303.549816
$ws.Range("O9").Value = 0.4120923552595624
$ws.Range("P9").Value = 0.4120923552595624
$ws.Range("Q9").Value = 458.4770551114026
$ws.Range("R9").Value = 4126.293496002623
$ws.Range("S9").Value = 0.1806977174509869
$ws.Range("T9").Value = 0.1806977174509869
# Row 10
$ws.Range("G10").Value = 3.895605666666667
$ws.Range("H10").Value = 11.686817
$ws.Range("I10").Value = 0.3769851220961256
$ws.Range("J10").Value = 0.3769851220961256
$ws.Range("M10").Value = 16.27546433333333
$ws.Range("N10").Value = 48.826393
$ws.Range("O10").Value = 0.06628560529319844
$ws.Range("P10").Value = 0.06628560529319844
$ws.Range("Q10").Value = 63.40279108456456
$ws.Range("R10").Value = 570.625119761081
$ws.Range("S10").Value = 0.024988687004672
$ws.Range("T10").Value = 0.024988687004672
# Row 11
$ws.Range("G11").Value = 3.895605666666667
$ws.Range("H11").Value = 11.686817
$ws.Range("I11").Value = 0.3769851220961256
$ws.Range("J11").Value = 0.3769851220961256
$ws.Range("O11").Value = 0.3480686258826592
$ws.Range("P11").Value = 0.3480686258826592
$ws.Range("Q11").Value = 332.9308418066773
$ws.Range("R11").Value = 2996.377576260096
$ws.Range("S11").Value = 0.131216693426205
$ws.Range("T11").Value = 0.131216693426205
# Row 12
$ws.Range("G12").Value = 3.895605666666667
$ws.Range("H12").Value = 11.686817
$ws.Range("I12").Value = 0.3769851220961256
$ws.Range("J12").Value = 0.3769851220961256
$ws.Range("M12").Value = 42.61351133333333
$ws.Range("N12").Value = 127.840534
$ws.Range("O12").Value = 0.17355341356458
$ws.Range("P12").Value = 0.17355341356458
$ws.Range("Q12").Value = 166.0054362266976
$ws.Range("R12").Value = 1494.048926040278
$ws.Range("S12").Value = 0.06542705480284257
$ws.Range("T12").Value = 0.06542705480284258
# Row 13
$ws.Range("G13").Value = 3.895605666666667
$ws.Range("H13").Value = 11.686817
$ws.Range("I13").Value = 0.3769851220961256
$ws.Range("J13").Value = 0.3769851220961256
$ws.Range("M13").Value = 101.183272
$ws.Range("N13").Value = 303.549816
$ws.Range("O13").Value = 0.4120923552595624
$ws.Range("P13").Value = 0.4120923552595624
$ws.Range("Q13").Value = 394.1701277750747
$ws.Range("R13").Value = 3547.531149975672
$ws.Range("S13").Value = 0.1553526868624061
$ws.Range("T13").Value = 0.1553526868624061
# Row 14
$ws.Range("G14").Value = 0.1037266666666667
$ws.Range("H14").Value = 0.31118
$ws.Range("I14").Value = 0.01003782555112075
$ws.Range("J14").Value = 0.01003782555112075
$ws.Range("M14").Value = 16.27546433333333
$ws.Range("N14").Value = 48.826393
$ws.Range("O14").Value = 0.06628560529319844
$ws.Range("P14").Value = 0.06628560529319844
$ws.Range("Q14").Value = 1.688199663748889
$ws.Range("R14").Value = 15.19379697374
$ws.Range("S14").Value = 0.0006653633424835722
$ws.Range("T14").Value = 0.0006653633424835721
# Row 15
$ws.Range("G15").Value = 0.1037266666666667
$ws.Range("H15").Value = 0.31118
$ws.Range("I15").Value = 0.01003782555112075
$ws.Range("J15").Value = 0.01003782555112075
$ws.Range("O15").Value = 0.3480686258826592
$ws.Range("P15").Value = 0.3480686258826592
$ws.Range("Q15").Value = 8.864810611255558
$ws.Range("R15").Value = 79.78329550130002
$ws.Range("S15").Value = 0.003493852146428446
$ws.Range("T15").Value = 0.003493852146428446
# Row 16
$ws.Range("G16").Value = 0.1037266666666667
$ws.Range("H16").Value = 0.31118
$ws.Range("I16").Value = 0.01003782555112075
$ws.Range("J16").Value = 0.01003782555112075
$ws.Range("M16").Value = 42.61351133333333
$ws.Range("N16").Value = 127.840534
$ws.Range("O16").Value = 0.17355341356458
$ws.Range("P16").Value = 0.17355341356458
$ws.Range("Q16").Value = 4.420157485568889
$ws.Range("R16").Value = 39.78141737012
$ws.Range("S16").Value = 0.001742098889162768
$ws.Range("T16").Value = 0.001742098889162768
# Row 17
$ws.Range("G17").Value = 0.1037266666666667
$ws.Range("H17").Value = 0.31118
$ws.Range("I17").Value = 0.01003782555112075
$ws.Range("J17").Value = 0.01003782555112075
$ws.Range("M17").Value = 101.183272
$ws.Range("N17").Value = 303.549816
$ws.Range("O17").Value = 0.4120923552595624
$ws.Range("P17").Value = 0.4120923552595624
$ws.Range("Q17").Value = 10.49540352698667
$ws.Range("R17").Value = 94.45863174287999
$ws.Range("S17").Value = 0.004136511173045965
$ws.Range("T17").Value = 0.004136511173045965
